# This script reproduces the change described in the commit "Fruta / hortaliza, semanal".
# A new weekly price-report row is inserted right after the existing row 60
# (i.e. as the new row 61), pushing all the following rows (old 61..98) down
# by one (new 62..99). The worksheet dimension grows from A1:T98 to A1:T99.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 61; this shifts rows 61-98 down to 62-99
# and automatically extends the sheet dimension to A1:T99.
$ws.Rows("61:61").Insert()

# Populate the newly inserted row 61 with the new weekly record.
$ws.Cells.Item(61, 1).Value  = 4
$ws.Cells.Item(61, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(61, 3).Value  = "Los Lagos"
$ws.Cells.Item(61, 4).Value  = 44777
$ws.Cells.Item(61, 5).Value  = 10
$ws.Cells.Item(61, 6).Value  = "Fruta"
$ws.Cells.Item(61, 7).Value  = 100104
$ws.Cells.Item(61, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(61, 9).Value  = 100104003
$ws.Cells.Item(61, 10).Value = "Membrillo"
$ws.Cells.Item(61, 11).Value = "Champion"
$ws.Cells.Item(61, 12).Value = "Primera"
$ws.Cells.Item(61, 13).Value = 300
$ws.Cells.Item(61, 14).Value = 13000
$ws.Cells.Item(61, 15).Value = 14000
$ws.Cells.Item(61, 16).Value = 13500
$ws.Cells.Item(61, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(61, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(61, 19).Value = 750
$ws.Cells.Item(61, 20).Value = 18

# Keep the date-formatted style (numFmt 165) on the new D61 cell, matching
# the rest of the "Fecha" column.
$ws.Cells.Item(61, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
